$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "29.420.14"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.848.68"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "240.40"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.6271"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.56%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.07694"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.28%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.2911"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.39%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "24.76"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.24%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07751"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.850.72"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.54%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.025"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("E14").Value = "  +3.70%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6811"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.33%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "83.56"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "6.173"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "29.452.76"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "228.55"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.38"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.422"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "157.21"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.1374"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.22%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "8.403"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "17.70"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.76%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.344"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.80%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.463"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.55%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.05644"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.119"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.032"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.841"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.162"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.43%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.7083"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.49%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.592"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.228.76"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.31%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.768"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01788"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.08%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.449"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.59%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.9057"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.48%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.033.05"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "101.71"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "65.85"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.21%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "7.166"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.00000000118"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "

$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.4010"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.1157"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.98%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "8.990"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.674"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
